# Update Leve profit-tracking figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the per-job Leve sheets, per the scheduled Sheets runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End / Animal Glue
$ws.Range("H5").Value = 340.9
$ws.Range("I5").Value = 377
$ws.Range("J5").Value = 16
$ws.Range("K5").Value = 377
$ws.Range("L5").Value = 16
$ws.Range("M5").Value = -262
$ws.Range("N5").Value = -246

# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 1831.2142
$ws.Range("I19").Value = 372.5
$ws.Range("K19").Value = 372.5
$ws.Range("M19").Value = -197.5

# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 1515.5883
$ws.Range("I28").Value = 1080.2222
$ws.Range("J28").Value = 2005.375
$ws.Range("K28").Value = 1080.2222
$ws.Range("L28").Value = 2005.375
$ws.Range("M28").Value = -595.2221999999999
$ws.Range("N28").Value = -2975.375

# Row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 3571.238
$ws.Range("J32").Value = 4086.8667
$ws.Range("L32").Value = 4086.8667
$ws.Range("N32").Value = -4738.8667

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 482.72726
$ws.Range("I98").Value = 419
$ws.Range("K98").Value = 419
$ws.Range("M98").Value = 1079

# Row 115: 5-bell Energy / Competent Craftsman's Syrup
$ws.Range("H115").Value = 485.66666
$ws.Range("I115").Value = 518.2
$ws.Range("J115").Value = 323
$ws.Range("K115").Value = 1554.6
$ws.Range("L115").Value = 969
$ws.Range("M115").Value = 12.39999999999986
$ws.Range("N115").Value = -4103

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 482.72726
$ws.Range("I122").Value = 419
$ws.Range("K122").Value = 1257
$ws.Range("M122").Value = 1193

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 1566.9296
$ws.Range("I132").Value = 1553.1177
$ws.Range("K132").Value = 4659.3531
$ws.Range("M132").Value = -2129.3531

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 9112.375
$ws.Range("I137").Value = 6185.643
$ws.Range("J137").Value = 13209.8
$ws.Range("K137").Value = 18556.929
$ws.Range("L137").Value = 39629.39999999999
$ws.Range("M137").Value = -16006.929
$ws.Range("N137").Value = -44729.39999999999

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 52633596
$ws.Range("I138").Value = 1518.375
$ws.Range("K138").Value = 4555.125
$ws.Range("M138").Value = 584.875

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 974.91113
$ws.Range("I141").Value = 841.525
$ws.Range("J141").Value = 2042
$ws.Range("K141").Value = 2524.575
$ws.Range("L141").Value = 6126
$ws.Range("M141").Value = 2655.425
$ws.Range("N141").Value = -16486

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3764.014
$ws.Range("I32").Value = 3872.7742
$ws.Range("J32").Value = 3089.7
$ws.Range("K32").Value = 3872.7742
$ws.Range("L32").Value = 3089.7
$ws.Range("M32").Value = -3585.7742
$ws.Range("N32").Value = -3663.7

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4947.375
$ws.Range("I61").Value = 3635.8572
$ws.Range("K61").Value = 3635.8572
$ws.Range("M61").Value = -3423.8572

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 40663.55
$ws.Range("I74").Value = 40663.55
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 40663.55
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -39789.55
$ws.Range("N74").ClearContents()

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 40663.55
$ws.Range("I77").Value = 40663.55
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 203317.75
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -198949.75
$ws.Range("N77").ClearContents()

# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 1600.5
$ws.Range("I97").Value = 1467.3334
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1467.3334
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -971.3334
$ws.Range("N97").Value = -2992

# Row 124: Ace of Gloves / High Durium Gauntlets of Fending
$ws.Range("H124").Value = 31858.143
$ws.Range("J124").Value = 36334.5
$ws.Range("L124").Value = 36334.5
$ws.Range("N124").Value = -46154.5

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 117030.43
$ws.Range("I132").Value = 2608.1667
$ws.Range("J132").Value = 269593.44
$ws.Range("K132").Value = 7824.500100000001
$ws.Range("L132").Value = 808780.3200000001
$ws.Range("M132").Value = -5294.500100000001
$ws.Range("N132").Value = -813840.3200000001

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4947.375
$ws.Range("I136").Value = 3635.8572
$ws.Range("K136").Value = 10907.5716
$ws.Range("M136").Value = -8357.571599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 75: I Saw the Pine / Hardsilver Saw
$ws.Range("H75").Value = 54406.11
$ws.Range("I75").Value = 9883.799999999999
$ws.Range("K75").Value = 9883.799999999999
$ws.Range("M75").Value = -8947.799999999999

# Row 78: I Came, I Sawed, I Conquered (L) / Hardsilver Saw
$ws.Range("H78").Value = 54406.11
$ws.Range("I78").Value = 9883.799999999999
$ws.Range("K78").Value = 29651.4
$ws.Range("M78").Value = -24971.4

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 24653
$ws.Range("I86").Value = 15824.8
$ws.Range("K86").Value = 15824.8
$ws.Range("M86").Value = -14701.8

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 24653
$ws.Range("I89").Value = 15824.8
$ws.Range("K89").Value = 79124
$ws.Range("M89").Value = -73508

# Row 96: Hammer Time / High Steel Sledgehammer
$ws.Range("H96").Value = 79424.25
$ws.Range("I96").Value = 9848.75
$ws.Range("K96").Value = 9848.75
$ws.Range("M96").Value = -7102.75

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1679.8572
$ws.Range("I107").Value = 1280.0303
$ws.Range("K107").Value = 1280.0303
$ws.Range("M107").Value = 639.9697000000001

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2338.111
$ws.Range("I134").Value = 1792.35
$ws.Range("J134").Value = 3897.4285
$ws.Range("K134").Value = 5377.049999999999
$ws.Range("L134").Value = 11692.2855
$ws.Range("M134").Value = -2842.049999999999
$ws.Range("N134").Value = -16762.2855

$ws = $wb.Worksheets.Item("CRP")
# Row 11: Leaving without Leave / Bronze Spear
$ws.Range("H11").Value = 2586.6667
$ws.Range("I11").Value = 1505
$ws.Range("J11").Value = 4750
$ws.Range("K11").Value = 1505
$ws.Range("L11").Value = 4750
$ws.Range("M11").Value = -1365
$ws.Range("N11").Value = -5030

# Row 21: Nightmare on My Street / Elm Cane
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 6631.552
$ws.Range("I31").Value = 3907.1428
$ws.Range("J31").Value = 9174.333000000001
$ws.Range("K31").Value = 3907.1428
$ws.Range("L31").Value = 9174.333000000001
$ws.Range("M31").Value = -3612.1428
$ws.Range("N31").Value = -9764.333000000001

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 6631.552
$ws.Range("I34").Value = 3907.1428
$ws.Range("J34").Value = 9174.333000000001
$ws.Range("K34").Value = 3907.1428
$ws.Range("L34").Value = 9174.333000000001
$ws.Range("M34").Value = -3705.1428
$ws.Range("N34").Value = -9578.333000000001

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 3224.4644
$ws.Range("I58").Value = 2612.1667
$ws.Range("K58").Value = 2612.1667
$ws.Range("M58").Value = -2409.1667

# Row 59: Bow Down to Magic / Crab Bow
$ws.Range("H59").Value = 42996.332
$ws.Range("J59").Value = 42999
$ws.Range("L59").Value = 42999
$ws.Range("N59").Value = -45289

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 5812.4
$ws.Range("I99").Value = 5812.4
$ws.Range("K99").Value = 5812.4
$ws.Range("M99").Value = -4314.4

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 5605.15
$ws.Range("I105").Value = 6736.0713
$ws.Range("J105").Value = 2966.3333
$ws.Range("K105").Value = 6736.0713
$ws.Range("L105").Value = 2966.3333
$ws.Range("M105").Value = -4989.0713
$ws.Range("N105").Value = -6460.3333

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 663.8125
$ws.Range("I107").Value = 424.66666
$ws.Range("J107").Value = 807.3
$ws.Range("K107").Value = 424.66666
$ws.Range("L107").Value = 807.3
$ws.Range("M107").Value = 1495.33334
$ws.Range("N107").Value = -4647.3

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 2987.125
$ws.Range("I122").Value = 1524.5
$ws.Range("J122").Value = 4449.75
$ws.Range("K122").Value = 4573.5
$ws.Range("L122").Value = 13349.25
$ws.Range("M122").Value = -2123.5
$ws.Range("N122").Value = -18249.25

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 5812.4
$ws.Range("I126").Value = 5812.4
$ws.Range("K126").Value = 17437.2
$ws.Range("M126").Value = -14967.2

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2191.5151
$ws.Range("I132").Value = 2234.5557
$ws.Range("J132").Value = 1997.8334
$ws.Range("K132").Value = 6703.6671
$ws.Range("L132").Value = 5993.5002
$ws.Range("M132").Value = -4173.6671
$ws.Range("N132").Value = -11053.5002

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 3272.96
$ws.Range("I134").Value = 2534.5264
$ws.Range("J134").Value = 5611.3335
$ws.Range("K134").Value = 7603.5792
$ws.Range("L134").Value = 16834.0005
$ws.Range("M134").Value = -5068.5792
$ws.Range("N134").Value = -21904.0005

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 3224.4644
$ws.Range("I136").Value = 2612.1667
$ws.Range("K136").Value = 7836.500100000001
$ws.Range("M136").Value = -5286.500100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 2577.4285
$ws.Range("I5").Value = 1329.3334
$ws.Range("K5").Value = 3988.0002
$ws.Range("M5").Value = -3876.0002

# Row 14: Keep Your Powder Dry / Kukuru Powder
$ws.Range("H14").Value = 296.46155
$ws.Range("I14").Value = 296.46155
$ws.Range("K14").Value = 889.38465
$ws.Range("M14").Value = -716.38465

# Row 58: Bread in the Clouds / La Noscean Toast
$ws.Range("H58").Value = 2834.3333
$ws.Range("I58").Value = 3752
$ws.Range("J58").Value = 999
$ws.Range("K58").Value = 11256
$ws.Range("L58").Value = 2997
$ws.Range("M58").Value = -11128
$ws.Range("N58").Value = -3253

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 418.74075
$ws.Range("J107").Value = 406
$ws.Range("L107").Value = 1218
$ws.Range("N107").Value = -5058

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 2218.4666
$ws.Range("I113").Value = 388
$ws.Range("J113").Value = 2676.0833
$ws.Range("K113").Value = 1164
$ws.Range("L113").Value = 8028.249899999999
$ws.Range("M113").Value = 1006
$ws.Range("N113").Value = -12368.2499

# Row 117: A Good Omen / Peppered Popotoes
$ws.Range("H117").Value = 1960.3889
$ws.Range("J117").Value = 3059.111
$ws.Range("L117").Value = 9177.332999999999
$ws.Range("N117").Value = -16061.333

# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 1442.6666
$ws.Range("I129").Value = 726.7143
$ws.Range("K129").Value = 2180.1429
$ws.Range("M129").Value = 2819.8571

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1499.1904
$ws.Range("J131").Value = 1675.4286
$ws.Range("L131").Value = 5026.2858
$ws.Range("N131").Value = -15106.2858

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 2577.4285
$ws.Range("I135").Value = 1329.3334
$ws.Range("K135").Value = 11964.0006
$ws.Range("M135").Value = -9429.000599999999

# Row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws.Range("H140").Value = 870.1667
$ws.Range("I140").Value = 870.1667
$ws.Range("K140").Value = 2610.5001
$ws.Range("M140").Value = 2569.4999

$ws = $wb.Worksheets.Item("GSM")
# Row 3: Needful Rings / Copper Wristlets
$ws.Range("H3").Value = 699.7143
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -384

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1451.5
$ws.Range("I102").Value = 1444.5714
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1444.5714
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 177.4286
$ws.Range("N102").Value = -4744

# Row 108: Satisfactory Sewing / Stonegold Needle
$ws.Range("H108").Value = 130000
$ws.Range("J108").Value = 130000
$ws.Range("L108").Value = 130000
$ws.Range("N108").Value = -137680

# Row 130: Planisphere to Paper / Chondrite Magitek Planisphere
$ws.Range("H130").Value = 56664.332
$ws.Range("J130").Value = 56664.332
$ws.Range("L130").Value = 56664.332
$ws.Range("N130").Value = -66704.33199999999

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 848.7646999999999
$ws.Range("I132").Value = 842.3333
$ws.Range("J132").Value = 897
$ws.Range("K132").Value = 2526.9999
$ws.Range("L132").Value = 2691
$ws.Range("M132").Value = 3.000100000000202
$ws.Range("N132").Value = -7751

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 3545.1667
$ws.Range("I7").Value = 1221.75
$ws.Range("K7").Value = 1221.75
$ws.Range("M7").Value = -1109.75

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 2295.3635
$ws.Range("I46").Value = 824.9286
$ws.Range("K46").Value = 824.9286
$ws.Range("M46").Value = -636.9286

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1758.5333
$ws.Range("J93").Value = 2233
$ws.Range("L93").Value = 2233
$ws.Range("N93").Value = -4729

# Row 94: Fitting In / Gaganaskin Hat of Aiming
$ws.Range("H94").Value = 85000
$ws.Range("J94").Value = 85000
$ws.Range("L94").Value = 85000
$ws.Range("N94").Value = -86352

# Row 99: Shoe on the Other Foot / Tigerskin Boots of Crafting
$ws.Range("H99").Value = 38784.75
$ws.Range("J99").Value = 48095
$ws.Range("L99").Value = 48095
$ws.Range("N99").Value = -54085

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 4536.778
$ws.Range("I122").Value = 3018.9285
$ws.Range("K122").Value = 9056.7855
$ws.Range("M122").Value = -6606.7855

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 3545.1667
$ws.Range("I126").Value = 1221.75
$ws.Range("K126").Value = 3665.25
$ws.Range("M126").Value = -1195.25

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 7464.4585
$ws.Range("I132").Value = 3998.7896
$ws.Range("K132").Value = 11996.3688
$ws.Range("M132").Value = -9466.3688

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 2751.5625
$ws.Range("I136").Value = 2707.375
$ws.Range("J136").Value = 2884.125
$ws.Range("K136").Value = 8122.125
$ws.Range("L136").Value = 8652.375
$ws.Range("M136").Value = -5572.125
$ws.Range("N136").Value = -13752.375

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 664.8333
$ws.Range("I107").Value = 663.3200000000001
$ws.Range("K107").Value = 1989.96
$ws.Range("M107").Value = -69.96000000000004

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1841.0465
$ws.Range("I132").Value = 1497.697
$ws.Range("K132").Value = 4493.090999999999
$ws.Range("M132").Value = -1963.090999999999

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 8666.729499999999
$ws.Range("I136").Value = 7766.8276
$ws.Range("K136").Value = 23300.4828
$ws.Range("M136").Value = -20750.4828

